# Auto commit at 2026-02-04  8:49:04.56
# Append two more days (2026-02-02 and 2026-02-03) of per-station charging
# data to the "daydata" table on Sheet1, as rows 4-7, matching the layout
# and formatting of the existing rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting (styles/number formats) of the existing data
# rows (2-3) down into the new rows (4-5 and 6-7) before writing values, so
# the new cells pick up the same style indices (date format, currency
# format, integer format) instead of minting new ones.
$ws.Range("A2:F3").Copy($ws.Range("A4"))
$ws.Range("A2:F3").Copy($ws.Range("A6"))

# Row 4: 2026-02-02, 四方坪站
$ws.Cells.Item(4, 1).Value = 46055
$ws.Cells.Item(4, 2).Value = "四方坪站"
$ws.Cells.Item(4, 3).Value = 11393.21
$ws.Cells.Item(4, 4).Value = 10328.209999999999
$ws.Cells.Item(4, 5).Value = 4243.25
$ws.Cells.Item(4, 6).Value = 484

# Row 5: 2026-02-02, 高岭站
$ws.Cells.Item(5, 1).Value = 46055
$ws.Cells.Item(5, 2).Value = "高岭站"
$ws.Cells.Item(5, 3).Value = 4732.8100000000004
$ws.Cells.Item(5, 4).Value = 4292.97
$ws.Cells.Item(5, 5).Value = 1305.55
$ws.Cells.Item(5, 6).Value = 184

# Row 6: 2026-02-03, 四方坪站
$ws.Cells.Item(6, 1).Value = 46056
$ws.Cells.Item(6, 2).Value = "四方坪站"
$ws.Cells.Item(6, 3).Value = 9959.92
$ws.Cells.Item(6, 4).Value = 8975.68
$ws.Cells.Item(6, 5).Value = 3735.71
$ws.Cells.Item(6, 6).Value = 436

# Row 7: 2026-02-03, 高岭站
$ws.Cells.Item(7, 1).Value = 46056
$ws.Cells.Item(7, 2).Value = "高岭站"
$ws.Cells.Item(7, 3).Value = 4281.92
$ws.Cells.Item(7, 4).Value = 3850.72
$ws.Cells.Item(7, 5).Value = 1142.08
$ws.Cells.Item(7, 6).Value = 157

# Match the workbook's recorded selection after the edit.
$ws.Range("H10").Select()
